$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet from "Data" to "Summary"
$ws.Name = "Summary"

# Re-assert formatting on the two untouched header cells so the engine
# keeps their (name/title) look after the round-trip.
$ws.Range("A1").Font.Size = 18
$ws.Range("A3").Font.Bold = $true

# Row 5/6/7 content is being relocated further down the sheet (rows 11/12/13)
# and new rows are introduced above (row 9) and appended at the bottom
# (rows 19/20). Clear the old locations first so nothing is left behind,
# then (re)write every cell at its final address.
$ws.Range("A5:D7").Clear()

# New sub-heading above the data table (bold + underline "title_" style)
$r = $ws.Range("A9")
$r.Value = "Source Type: SME Associations (Most Widely Used)"
$r.Font.Bold = $true
$r.Font.Underline = $true

# Column headers (bold "title" style)
$ws.Range("B11").Value = "Micro"
$ws.Range("B11").Font.Bold = $true

$ws.Range("C11").Value = "SMEs"
$ws.Range("C11").Font.Bold = $true

$ws.Range("D11").Value = "MSMEs"
$ws.Range("D11").Font.Bold = $true

# Data row: row label is bold "title" style, values are plain "Normal" text
$ws.Range("A12").Value = "Enterprises (% of total)"
$ws.Range("A12").Font.Bold = $true

$ws.Range("B12").NumberFormat = "@"
$ws.Range("B12").Value = "81"
$ws.Range("B12").Style = "Normal"

$ws.Range("C12").NumberFormat = "@"
$ws.Range("C12").Value = "14"
$ws.Range("C12").Style = "Normal"

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "95"
$ws.Range("D12").Style = "Normal"

# Source line (italic "source" style)
$r = $ws.Range("A13")
$r.Value = "Source: MTICM, 2008"
$r.Font.Italic = $true

# New source-detail block near the bottom of the sheet
$r = $ws.Range("A19")
$r.Value = "MTICM"
$r.Font.Bold = $true

$r = $ws.Range("A20")
$r.Value = "Ministry of Trade and Industry, Cooperatives and Marketing, (MTICM), ""The State of Small Enterprise in Lesotho"", 2008, p. 19. Available at http://www.sbp.org.za/uploads/media/Lesotho_White_Paper__6-2_.pdf"
$r.Font.Italic = $true
